$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A width (cosmetic) ---
$ws.Columns("A:A").ColumnWidth = 13

# --- helper: a couple of cells that already carry the shared "Area" style (fontId 1 / s="1") ---
# H2 / C8 are pre-existing cells styled that way; we paste-special (formats only) from them
# onto any newly written cell so the new cell ends up on style index 1 as well, matching the
# original workbook's look without inflating styles.xml with duplicate style records.

function Set-AreaStyle($addr) {
    $ws.Range("H2").Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
}

function Set-MateriaStyle($addr) {
    $ws.Range("C8").Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
}

# ============ Existing rows 20-37: fill in weekly/total hours + area ============

# Row 20 - Analisis matematico I
$ws.Range("F20").Value = 6
$ws.Range("G20").Value = 96
$ws.Range("H20").Value = "CBGyE"
Set-AreaStyle "H20"

# Row 21 - Bases de Datos I
$ws.Range("F21").Value = 6
$ws.Range("G21").Value = 96
$ws.Range("H21").Value = "ISBDSI"
Set-AreaStyle "H21"

# Row 22 - Redes y Comunicaciones
$ws.Range("F22").Value = 6
$ws.Range("G22").Value = 96
$ws.Range("H22").Value = "ARSORE"
Set-AreaStyle "H22"

# Row 23 - Seminario de Justicia y Derechos Humanos
$ws.Range("F23").Value = 6
$ws.Range("G23").Value = 64
$ws.Range("H23").Value = "ASyP"
Set-AreaStyle "H23"

# Row 24 - Practicas Pre-profesionales I (note: no F value, and H stays unstyled like in the source)
$ws.Range("G24").Value = 80
$ws.Range("H24").Value = "Espacios de integración"

# Row 25 - Conceptos y Paradigmas de Lenguajes de Programacion
$ws.Range("F25").Value = 4
$ws.Range("G25").Value = 64
$ws.Range("H25").Value = "AyL"
Set-AreaStyle "H25"

# Row 26 - Orientacion a Objetos II
$ws.Range("F26").Value = 4
$ws.Range("G26").Value = 64
$ws.Range("H26").Value = "AyL"
Set-AreaStyle "H26"

# Row 27 - Ingenieria de Software III
$ws.Range("F27").Value = 6
$ws.Range("G27").Value = 96
$ws.Range("H27").Value = "ISBDSI"
Set-AreaStyle "H27"

# Row 28 - Bases de Datos II
$ws.Range("F28").Value = 6
$ws.Range("G28").Value = 96
$ws.Range("H28").Value = "ISBDSI"
Set-AreaStyle "H28"

# Row 29 - Sistemas Operativos
$ws.Range("F29").Value = 6
$ws.Range("G29").Value = 96
$ws.Range("H29").Value = "ARSORE"
Set-AreaStyle "H29"

# Row 30 - Analisis matematico II
$ws.Range("F30").Value = 6
$ws.Range("G30").Value = 96
$ws.Range("H30").Value = "CBGyE"
Set-AreaStyle "H30"

# Row 31 - Fundamentos de Teoria de la Computacion
$ws.Range("F31").Value = 6
$ws.Range("G31").Value = 96
$ws.Range("H31").Value = "CBGyE"
Set-AreaStyle "H31"

# Row 32 - Desarrollo de Software en Sistemas Distribuidos
$ws.Range("F32").Value = 6
$ws.Range("G32").Value = 96
$ws.Range("H32").Value = "ISBDSI"
Set-AreaStyle "H32"

# Row 33 - Sistemas y Organizaciones
$ws.Range("F33").Value = 6
$ws.Range("G33").Value = 96
$ws.Range("H33").Value = "ISBDSI"
Set-AreaStyle "H33"

# Row 34 - Proyecto de Software
$ws.Range("F34").Value = 6
$ws.Range("G34").Value = 96
$ws.Range("H34").Value = "ISBDSI"
Set-AreaStyle "H34"

# Row 35 - Taller de Metodologia de la Investigacion (correlativas updated too)
$ws.Range("D35").Value = '["Ingenieria de Software II","Bases de Datos I","Redes y Comunicaciones","Orientacion a Objetos II","Sistemas Operativos"]'
$ws.Range("F35").Value = 6
$ws.Range("G35").Value = 96
$ws.Range("H35").Value = "CBGyE"
Set-AreaStyle "H35"

# Row 36 - Probabilidad y Estadistica
$ws.Range("F36").Value = 6
$ws.Range("G36").Value = 96
$ws.Range("H36").Value = "CBGyE"
Set-AreaStyle "H36"

# Row 37 - Seminario Optativo* (no area assigned)
$ws.Range("F37").Value = 4
$ws.Range("G37").Value = 64

# ============ New rows 38-41: 10th-semester / plan-2025 subjects ============

# Row 38 - Escenarios Tecnologicos
$ws.Range("A38").Value = "Lic. Sistemas"
$ws.Range("B38").Value = 2024
$ws.Range("C38").Value = "Escenarios Tecnologicos "
Set-MateriaStyle "C38"
$ws.Range("E38").Value = 10
$ws.Range("F38").Value = 4
$ws.Range("G38").Value = 64
$ws.Range("H38").Value = "ASyP"
Set-AreaStyle "H38"

# Row 39 - Aspectos sociales y profesionales de la Informatica
$ws.Range("A39").Value = "Lic. Sistemas"
$ws.Range("B39").Value = 2024
$ws.Range("C39").Value = "Aspectos sociales y profesionales de la Informatica"
Set-MateriaStyle "C39"
$ws.Range("D39").Value = '["Sistemas y Organizaciones"]'
$ws.Range("E39").Value = 10
$ws.Range("F39").Value = 4
$ws.Range("G39").Value = 64
$ws.Range("H39").Value = "ASyP"
Set-AreaStyle "H39"

# Row 40 - Practicas Preprofesionales II
$ws.Range("A40").Value = "Lic. Sistemas"
$ws.Range("B40").Value = 2024
$ws.Range("C40").Value = "Practicas Preprofesionales II"
Set-MateriaStyle "C40"
$ws.Range("D40").Value = '["Practicas Pre-profesionales I"]'
$ws.Range("E40").Value = 10
$ws.Range("G40").Value = 7
$ws.Range("H40").Value = "Espacios de integración"
Set-AreaStyle "H40"

# Row 41 - Taller de Proyectos I+D+i
$ws.Range("A41").Value = "Lic. Sistemas"
$ws.Range("B41").Value = 2024
$ws.Range("C41").Value = "Taller de Proyectos I+D+i"
Set-MateriaStyle "C41"
$ws.Range("D41").Value = '["Proyecto de Software","Taller de Metodologia de la Investigacion","Probabilidad y Estadistica"]'
$ws.Range("E41").Value = 10
$ws.Range("F41").Value = 6
$ws.Range("G41").Value = 64
$ws.Range("H41").Value = "Espacios de integración"
Set-AreaStyle "H41"

# --- Page setup: A4 portrait (matches the updated pageSetup element) ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection / view state, best effort ---
$ws.Range("H41").Select()

Write-Host "Edit applied"
